$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "43.140.30"
$ws.Range("E2").Value = "  +4.93%  "
$ws.Range("D3").Value = "2.241.30"
$ws.Range("E3").Value = "  +4.89%  "
$ws.Range("E4").Value = "  -0.17%  "
$ws.Range("D5").Value = "'251.93"
$ws.Range("E5").Value = "  +7.58%  "
$ws.Range("E6").Value = "  +2.88%  "
$ws.Range("D7").Value = "'75.31"
$ws.Range("E7").Value = "  +9.51%  "
$ws.Range("E8").Value = "  -0.22%  "
$ws.Range("D9").Value = "'0.601"
$ws.Range("E9").Value = "  +6.99%  "
$ws.Range("D10").Value = "'41.37"
$ws.Range("E10").Value = "  +8.65%  "
$ws.Range("D11").Value = "'0.0930"
$ws.Range("E11").Value = "  +4.71%  "
$ws.Range("E12").Value = "  +5.72%  "
$ws.Range("D13").Value = "'0.102"
$ws.Range("E13").Value = "  +2.54%  "
$ws.Range("D14").Value = "2.577.31"
$ws.Range("E14").Value = "  +4.86%  "
$ws.Range("D15").Value = "'14.66"
$ws.Range("E15").Value = "  +2.93%  "
$ws.Range("D16").Value = "2.249.30"
$ws.Range("E16").Value = "  +5.51%  "
$ws.Range("D17").Value = "'0.793"
$ws.Range("E17").Value = "  +2.75%  "
$ws.Range("D18").Value = "43.032.87"
$ws.Range("E18").Value = "  +4.99%  "
$ws.Range("E19").Value = "  +6.39%  "
$ws.Range("D20").Value = "'71.17"
$ws.Range("E20").Value = "  +3.52%  "
$ws.Range("E21").Value = "  +5.47%  "
$ws.Range("D22").Value = "'228.49"
$ws.Range("E22").Value = "  +2.21%  "
$ws.Range("E23").Value = "  +16.39%  "
$ws.Range("E24").Value = "  +2.34%  "
$ws.Range("E25").Value = "  +0.02%  "
$ws.Range("D26").Value = "'10.76"
$ws.Range("E26").Value = "  +2.96%  "
$ws.Range("E27").Value = "  +2.82%  "
$ws.Range("D28").Value = "'39.47"
$ws.Range("E28").Value = "  +29.49%  "
$ws.Range("E29").Value = "  +6.33%  "
$ws.Range("D30").Value = "'2.24"
$ws.Range("E30").Value = "  +4.38%  "
$ws.Range("D31").Value = "'171.93"
$ws.Range("E31").Value = "  +1.90%  "
$ws.Range("D32").Value = "'20.24"
$ws.Range("E32").Value = "  +4.06%  "
$ws.Range("D33").Value = "'0.0802"
$ws.Range("E33").Value = "  +7.75%  "
$ws.Range("E35").Value = "  +2.69%  "
$ws.Range("E36").Value = "  +10.93%  "
$ws.Range("D37").Value = "'4.50"
$ws.Range("E37").Value = "  +11.01%  "
$ws.Range("D38").Value = "'0.0331"
$ws.Range("E38").Value = "  +18.84%  "
$ws.Range("D39").Value = "'12.96"
$ws.Range("E39").Value = "  +12.71%  "
$ws.Range("E40").Value = "  +4.58%  "
$ws.Range("D41").Value = "'0.204"
$ws.Range("E41").Value = "  +11.05%  "
$ws.Range("D42").Value = "'5.44"
$ws.Range("E42").Value = "  +4.63%  "
$ws.Range("D43").Value = "'59.61"
$ws.Range("E43").Value = "  +4.94%  "
$ws.Range("D44").Value = "'8.67"
$ws.Range("E44").Value = "  +6.70%  "
$ws.Range("D45").Value = "'103.83"
$ws.Range("D46").Value = "'0.483"
$ws.Range("E46").Value = "  +32.42%  "
$ws.Range("D47").Value = "'0.0989"
$ws.Range("E47").Value = "  +4.46%  "
$ws.Range("E48").Value = "  +14.77%  "
$ws.Range("E49").Value = "  +4.39%  "
$ws.Range("E50").Value = "  +5.30%  "
$ws.Range("E51").Value = "  +3.28%  "
